# The "Recorded By" column (G) stores comma-separated lists of the
# users/processes that recorded a session. A handful of rows still had
# these names in their original (legacy) order; this normalizes them by
# rotating the list left by one position (moving the first entry to the
# end), matching how the other rows in the sheet are already ordered.
#
# Only cells whose value is EXACTLY one of the three known legacy
# orderings are touched - everything else (including cells that already
# use the new ordering) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rotations = @{
    "dnasr281@gmail.com, System"            = "System, dnasr281@gmail.com"
    "admin@admin.com, System"               = "System, admin@admin.com"
    "system, backup@backdoor.com, System"   = "backup@backdoor.com, System, system"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }

    if ($rotations.ContainsKey($val)) {
        $cell.Value = $rotations[$val]
    }
}
